# Edit replicates the following source change:
#  1. Slide 5's table (graphicFrame "Google Shape;122;p17") switches its
#     table style from {126D691A-9FBF-4476-AEBA-CC6B0AB18413} to
#     {E2B98C09-6FE3-4A15-8C0C-74DA049C3844}.
#  2. The deck's theme color scheme is changed from the custom
#     "Integral" / "Red Violet" palette to the stock "Office" palette
#     (dk1/lt1 are unchanged; dk2, lt2 and all six accents + hyperlink
#     colors move to the Office defaults).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{E2B98C09-6FE3-4A15-8C0C-74DA049C3844}")

# --- 2. Swap the theme colors over to the stock Office palette -----------
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(3).RGB  = 6968388    # dk2      454551 -> 44546A
$colors.Item(4).RGB  = 15132391   # lt2      D8D9DC -> E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  E32D91 -> 5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  C830CC -> ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  4EA6DC -> A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  4775E7 -> FFC000
$colors.Item(9).RGB  = 12874308   # accent5  8971E1 -> 4472C4
$colors.Item(10).RGB = 4697456    # accent6  D54773 -> 70AD47
$colors.Item(11).RGB = 12673797   # hlink    6B9F25 -> 0563C1
$colors.Item(12).RGB = 7491477    # folHlink 8C8C8C -> 954F72
